$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.601.70"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3
$ws.Range("D3").Value = "2.466.64"
$ws.Range("E3").Value = "  -2.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.02"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.43"
$ws.Range("E6").Value = "  -4.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  -1.59%  "

# Row 9
$ws.Range("D9").Value = "2.477.83"
$ws.Range("E9").Value = "  -2.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  -0.38%  "

# Row 11
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  -1.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -3.81%  "

# Row 14
$ws.Range("D14").Value = "2.903.07"
$ws.Range("E14").Value = "  -2.99%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.92"
$ws.Range("E15").Value = "  -2.59%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "58.492.19"
$ws.Range("E16").Value = "  -1.87%  "

# Row 17
$ws.Range("E17").Value = "  -2.76%  "

# Row 18
$ws.Range("D18").Value = "2.476.18"
$ws.Range("E18").Value = "  -2.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.87"
$ws.Range("E19").Value = "  -2.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.21"
$ws.Range("E20").Value = "  -1.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.90"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.76"
$ws.Range("E23").Value = "  -1.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.74"
$ws.Range("E24").Value = "  -0.44%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.414"
$ws.Range("E25").Value = "  -2.76%  "

# Row 26
$ws.Range("E26").Value = "  -0.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.986"
$ws.Range("E27").Value = "  -1.30%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  -5.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.70"
$ws.Range("E29").Value = "  -3.38%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0758"
$ws.Range("E30").Value = "  -2.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.76"
$ws.Range("E31").Value = "  -2.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.66"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("E34").Value = "  -1.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.33"
$ws.Range("E35").Value = "  -1.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.36"
$ws.Range("E36").Value = "  -5.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.05"
$ws.Range("E37").Value = "  -6.34%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("E38").Value = "  -3.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.49"
$ws.Range("E39").Value = "  -1.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("E40").Value = "  -1.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("E41").Value = "  -2.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.17"
$ws.Range("E42").Value = "  -5.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.19"
$ws.Range("E43").Value = "  -5.70%  "

# Row 44
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.82"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.592"
$ws.Range("E46").Value = "  -1.82%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.50"
$ws.Range("E47").Value = "  -1.87%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0925"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").Value = "  -1.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0219"
$ws.Range("E50").Value = "  -3.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.30"
$ws.Range("E51").Value = "  -3.07%  "
